$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refresh Coin / Link / Price / Volume(1h) figures to match the latest coinranking.com scrape.
# Rows 38-39: RenderToken and Bittensor swap ranking positions with fresh data.

$ws.Range("D2").Value2 = '61.416.47'
$ws.Range("E2").Value2 = '  -2.50%  '

$ws.Range("D3").Value2 = '2.556.06'
$ws.Range("E3").Value2 = '  -4.49%  '

$ws.Range("E4").Value2 = '  +0.00%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value2 = '546.23'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value2 = '  -0.37%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value2 = '150.98'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value2 = '  -3.84%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value2 = '0.999'
$ws.Range("D7").Style = "Normal"

$ws.Range("E8").Value2 = '  -0.21%  '

$ws.Range("E9").Value2 = '  -2.40%  '

$ws.Range("E10").Value2 = '  -0.72%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value2 = '5.46'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value2 = '  +7.49%  '

$ws.Range("E12").Value2 = '  -1.35%  '

$ws.Range("D13").Value2 = '3.010.67'
$ws.Range("E13").Value2 = '  -4.27%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value2 = '25.09'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value2 = '  -3.69%  '

$ws.Range("D15").Value2 = '61.347.91'
$ws.Range("E15").Value2 = '  -2.43%  '

$ws.Range("E16").Value2 = '  -1.28%  '

$ws.Range("D17").Value2 = '2.567.07'
$ws.Range("E17").Value2 = '  -4.11%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value2 = '11.47'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value2 = '  -3.62%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value2 = '4.48'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value2 = '  -1.50%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value2 = '335.59'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value2 = '  -1.57%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value2 = '0.998'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value2 = '  +0.32%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value2 = '5.95'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value2 = '  -5.09%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value2 = '0.483'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value2 = '  -4.04%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value2 = '62.68'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value2 = '  -1.12%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value2 = '0.166'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value2 = '  -1.21%  '

$ws.Range("E26").Value2 = '  -0.02%  '

$ws.Range("E27").Value2 = '  -0.57%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value2 = '7.09'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value2 = '  +1.84%  '

$ws.Range("D29").Value2 = '0.0₃0813'
$ws.Range("E29").Value2 = '  -4.04%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value2 = '1.31'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value2 = '  -1.24%  '

$ws.Range("E31").Value2 = '  -2.61%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value2 = '160.76'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value2 = '  -2.99%  '

$ws.Range("E33").Value2 = '  +0.02%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value2 = '4.80'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value2 = '  +0.45%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value2 = '18.93'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value2 = '  -2.67%  '

$ws.Range("E36").Value2 = '  -2.84%  '

$ws.Range("E37").Value2 = '  -0.08%  '

$ws.Range("B38").Value2 = 'Bittensor'
$ws.Range("C38").Value2 = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value2 = '322.42'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value2 = '  -4.47%  '

$ws.Range("B39").Value2 = 'RenderToken'
$ws.Range("C39").Value2 = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value2 = '5.94'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value2 = '  -2.69%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value2 = '0.877'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value2 = '  -6.05%  '

$ws.Range("E41").Value2 = '  -1.10%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value2 = '37.17'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value2 = '  -2.31%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value2 = '0.999'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value2 = '  +0.08%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value2 = '20.38'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value2 = '  -1.33%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value2 = '10.92'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value2 = '  -1.12%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value2 = '0.601'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value2 = '  -2.38%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value2 = '0.0957'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value2 = '  -1.17%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value2 = '0.0537'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value2 = '  -3.80%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value2 = '19.23'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value2 = '  -4.52%  '

$ws.Range("E50").Value2 = '  -1.43%  '

$ws.Range("D51").Value2 = '2.034.97'
$ws.Range("E51").Value2 = '  -1.91%  '

